$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 152.28572
$ws.Range("I9").Value = 110.47369
$ws.Range("J9").Value = 549.5
$ws.Range("K9").Value = 110.47369
$ws.Range("L9").Value = 549.5
$ws.Range("M9").Value = 58.52631
$ws.Range("N9").Value = -887.5
$ws.Range("H64").Value = 5872.8096
$ws.Range("I64").Value = 5241.1113
$ws.Range("K64").Value = 5241.1113
$ws.Range("M64").Value = -4993.1113
$ws.Range("H67").Value = 5872.8096
$ws.Range("I67").Value = 5241.1113
$ws.Range("K67").Value = 5241.1113
$ws.Range("M67").Value = -4383.1113
$ws.Range("H76").Value = 6217
$ws.Range("I76").Value = 6739.8
$ws.Range("K76").Value = 6739.8
$ws.Range("M76").Value = -6424.8
$ws.Range("H79").Value = 6217
$ws.Range("I79").Value = 6739.8
$ws.Range("K79").Value = 6739.8
$ws.Range("M79").Value = -5647.8
$ws.Range("H103").Value = 31251296
$ws.Range("I103").Value = 740
$ws.Range("J103").Value = 45456092
$ws.Range("K103").Value = 2220
$ws.Range("L103").Value = 136368276
$ws.Range("M103").Value = -1634
$ws.Range("N103").Value = -136369448
$ws.Range("H113").Value = 7913.3335
$ws.Range("I113").Value = 7335
$ws.Range("J113").Value = 8491.666999999999
$ws.Range("K113").Value = 7335
$ws.Range("L113").Value = 8491.666999999999
$ws.Range("M113").Value = -4081
$ws.Range("N113").Value = -14999.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 648
$ws.Range("I2").Value = 631.1111
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 631.1111
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = -518.1111
$ws.Range("N2").Value = -1026
$ws.Range("H61").Value = 7245721.5
$ws.Range("I61").Value = 10005388
$ws.Range("J61").Value = 1113130.6
$ws.Range("K61").Value = 10005388
$ws.Range("L61").Value = 1113130.6
$ws.Range("M61").Value = -10005176
$ws.Range("N61").Value = -1113554.6
$ws.Range("H116").Value = 648
$ws.Range("I116").Value = 631.1111
$ws.Range("J116").Value = 800
$ws.Range("K116").Value = 631.1111
$ws.Range("L116").Value = 800
$ws.Range("M116").Value = 1662.8889
$ws.Range("N116").Value = -5388
$ws.Range("H136").Value = 7245721.5
$ws.Range("I136").Value = 10005388
$ws.Range("J136").Value = 1113130.6
$ws.Range("K136").Value = 30016164
$ws.Range("L136").Value = 3339391.8
$ws.Range("M136").Value = -30013614
$ws.Range("N136").Value = -3344491.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 648
$ws.Range("I3").Value = 631.1111
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 631.1111
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -517.1111
$ws.Range("N3").Value = -1028
$ws.Range("H20").Value = 2886.8684
$ws.Range("I20").Value = 2936.8
$ws.Range("K20").Value = 2936.8
$ws.Range("M20").Value = -2689.8
$ws.Range("H134").Value = 4765410
$ws.Range("I134").Value = 3277.6924
$ws.Range("J134").Value = 12503875
$ws.Range("K134").Value = 9833.0772
$ws.Range("L134").Value = 37511625
$ws.Range("M134").Value = -7298.0772
$ws.Range("N134").Value = -37516695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30327544
$ws.Range("I31").Value = 35742108
$ws.Range("J31").Value = 5970
$ws.Range("K31").Value = 35742108
$ws.Range("L31").Value = 5970
$ws.Range("M31").Value = -35741813
$ws.Range("N31").Value = -6560
$ws.Range("H34").Value = 30327544
$ws.Range("I34").Value = 35742108
$ws.Range("J34").Value = 5970
$ws.Range("K34").Value = 35742108
$ws.Range("L34").Value = 5970
$ws.Range("M34").Value = -35741906
$ws.Range("N34").Value = -6374
$ws.Range("H58").Value = 2504
$ws.Range("I58").Value = 2303.543
$ws.Range("K58").Value = 2303.543
$ws.Range("M58").Value = -2100.543
$ws.Range("H86").Value = 8812.166999999999
$ws.Range("I86").Value = 6319
$ws.Range("K86").Value = 6319
$ws.Range("M86").Value = -5196
$ws.Range("H89").Value = 8812.166999999999
$ws.Range("I89").Value = 6319
$ws.Range("K89").Value = 31595
$ws.Range("M89").Value = -25979
$ws.Range("H107").Value = 1958.1428
$ws.Range("I107").Value = 714.2857
$ws.Range("K107").Value = 714.2857
$ws.Range("M107").Value = 1205.7143
$ws.Range("H121").Value = 89375
$ws.Range("J121").Value = 89375
$ws.Range("L121").Value = 89375
$ws.Range("N121").Value = -91995
$ws.Range("H132").Value = 2803.3076
$ws.Range("I132").Value = 2677.0588
$ws.Range("J132").Value = 3041.7778
$ws.Range("K132").Value = 8031.176399999999
$ws.Range("L132").Value = 9125.3334
$ws.Range("M132").Value = -5501.176399999999
$ws.Range("N132").Value = -14185.3334
$ws.Range("H134").Value = 1396.36
$ws.Range("I134").Value = 1114.6842
$ws.Range("K134").Value = 3344.0526
$ws.Range("M134").Value = -809.0526
$ws.Range("H136").Value = 2504
$ws.Range("I136").Value = 2303.543
$ws.Range("K136").Value = 6910.629000000001
$ws.Range("M136").Value = -4360.629000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 15766
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 6000
$ws.Range("M22").Value = -5831
$ws.Range("H27").Value = 15766
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 6000
$ws.Range("M27").Value = -5898
$ws.Range("H38").Value = 66.666664
$ws.Range("H70").Value = 13969.857
$ws.Range("I70").Value = 6224.8
$ws.Range("K70").Value = 18674.4
$ws.Range("M70").Value = -18359.4
$ws.Range("H73").Value = 13969.857
$ws.Range("I73").Value = 6224.8
$ws.Range("K73").Value = 18674.4
$ws.Range("M73").Value = -17582.4
$ws.Range("H99").Value = 14785.667
$ws.Range("I99").Value = 25
$ws.Range("K99").Value = 75
$ws.Range("M99").Value = 2171
$ws.Range("H132").Value = 1982.7097
$ws.Range("I132").Value = 1251.7142
$ws.Range("J132").Value = 2195.9167
$ws.Range("K132").Value = 11265.4278
$ws.Range("L132").Value = 19763.2503
$ws.Range("M132").Value = -8735.427799999999
$ws.Range("N132").Value = -24823.2503
$ws.Range("H137").Value = 9252.933999999999
$ws.Range("I137").Value = 4766.1816
$ws.Range("J137").Value = 21591.5
$ws.Range("K137").Value = 14298.5448
$ws.Range("L137").Value = 64774.5
$ws.Range("M137").Value = -9198.5448
$ws.Range("N137").Value = -74974.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 10699.75
$ws.Range("J70").Value = 12052
$ws.Range("K70").Value = 10699.75
$ws.Range("L70").Value = 12052
$ws.Range("M70").Value = -10429.75
$ws.Range("N70").Value = -12592
$ws.Range("I73").Value = 10699.75
$ws.Range("J73").Value = 12052
$ws.Range("K73").Value = 10699.75
$ws.Range("L73").Value = 12052
$ws.Range("M73").Value = -9763.75
$ws.Range("N73").Value = -13924
$ws.Range("H122").Value = 5260.871
$ws.Range("I122").Value = 4771
$ws.Range("J122").Value = 5939.154
$ws.Range("K122").Value = 14313
$ws.Range("L122").Value = 17817.462
$ws.Range("M122").Value = -11863
$ws.Range("N122").Value = -22717.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 99999
$ws.Range("J18").Value = 99999
$ws.Range("L18").Value = 99999
$ws.Range("N18").Value = -100343
$ws.Range("H31").Value = 12507.5
$ws.Range("I31").Value = 5015
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 5015
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -4767
$ws.Range("N31").Value = -20496
$ws.Range("H35").Value = 2016.6666
$ws.Range("I35").Value = 2016.6666
$ws.Range("K35").Value = 2016.6666
$ws.Range("M35").Value = -1680.6666
$ws.Range("H100").Value = 20858582
$ws.Range("I100").Value = 3362.8572
$ws.Range("K100").Value = 3362.8572
$ws.Range("M100").Value = -2821.8572
$ws.Range("H122").Value = 4400.596
$ws.Range("I122").Value = 3461.6
$ws.Range("K122").Value = 10384.8
$ws.Range("M122").Value = -7934.799999999999
$ws.Range("H136").Value = 2521.9412
$ws.Range("I136").Value = 2219.1
$ws.Range("K136").Value = 6657.299999999999
$ws.Range("M136").Value = -4107.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 19999.8
$ws.Range("I18").Value = 20000
$ws.Range("J18").Value = 19999.5
$ws.Range("K18").Value = 20000
$ws.Range("L18").Value = 19999.5
$ws.Range("M18").Value = -19827
$ws.Range("N18").Value = -20345.5
$ws.Range("H46").Value = 128964.5
$ws.Range("J46").Value = 128964.5
$ws.Range("L46").Value = 128964.5
$ws.Range("N46").Value = -129426.5
$ws.Range("H113").Value = 606.61536
$ws.Range("I113").Value = 502
$ws.Range("J113").Value = 696.2857
$ws.Range("K113").Value = 1506
$ws.Range("L113").Value = 2088.8571
$ws.Range("M113").Value = 664
$ws.Range("N113").Value = -6428.8571
$ws.Range("H134").Value = 128964.5
$ws.Range("J134").Value = 128964.5
$ws.Range("L134").Value = 386893.5
$ws.Range("N134").Value = -391963.5
